{"js": "// Update the two-digit multiplication problems in the worksheet table.\n// Each \"AA\u00d7BB=\" expression is unique in the document, so a scoped,\n// case-sensitive search-and-replace for each old/new pair is safe.\nconst replacements = [\n  [\"62\u00d744=\", \"52\u00d725=\"],\n  [\"75\u00d736=\", \"56\u00d741=\"],\n  [\"28\u00d746=\", \"32\u00d725=\"],\n  [\"75\u00d755=\", \"96\u00d739=\"],\n  [\"52\u00d787=\", \"41\u00d728=\"],\n  [\"29\u00d760=\", \"13\u00d715=\"],\n  [\"37\u00d791=\", \"45\u00d737=\"],\n  [\"64\u00d796=\", \"74\u00d784=\"],\n  [\"21\u00d724=\", \"45\u00d761=\"],\n  [\"22\u00d748=\", \"24\u00d720=\"],\n  [\"92\u00d764=\", \"84\u00d792=\"],\n  [\"15\u00d719=\", \"15\u00d775=\"],\n  [\"88\u00d717=\", \"56\u00d750=\"],\n  [\"69\u00d748=\", \"36\u00d771=\"],\n  [\"61\u00d798=\", \"95\u00d771=\"],\n  [\"90\u00d798=\", \"99\u00d777=\"],\n  [\"47\u00d792=\", \"47\u00d768=\"],\n  [\"50\u00d725=\", \"28\u00d779=\"],\n  [\"35\u00d771=\", \"70\u00d741=\"],\n  [\"73\u00d721=\", \"20\u00d732=\"],\n  [\"95\u00d738=\", \"76\u00d786=\"],\n  [\"54\u00d782=\", \"56\u00d720=\"],\n  [\"25\u00d788=\", \"45\u00d740=\"],\n  [\"52\u00d788=\", \"16\u00d715=\"],\n  [\"93\u00d732=\", \"46\u00d761=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the two-digit multiplication problems in the worksheet table.\n# Each \"AA\u00d7BB=\" expression is unique in the document, so a straight\n# Find/Replace (wdReplaceAll) for each old/new pair is safe and idempotent.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"62\u00d744=\"; New = \"52\u00d725=\" },\n  @{ Old = \"75\u00d736=\"; New = \"56\u00d741=\" },\n  @{ Old = \"28\u00d746=\"; New = \"32\u00d725=\" },\n  @{ Old = \"75\u00d755=\"; New = \"96\u00d739=\" },\n  @{ Old = \"52\u00d787=\"; New = \"41\u00d728=\" },\n  @{ Old = \"29\u00d760=\"; New = \"13\u00d715=\" },\n  @{ Old = \"37\u00d791=\"; New = \"45\u00d737=\" },\n  @{ Old = \"64\u00d796=\"; New = \"74\u00d784=\" },\n  @{ Old = \"21\u00d724=\"; New = \"45\u00d761=\" },\n  @{ Old = \"22\u00d748=\"; New = \"24\u00d720=\" },\n  @{ Old = \"92\u00d764=\"; New = \"84\u00d792=\" },\n  @{ Old = \"15\u00d719=\"; New = \"15\u00d775=\" },\n  @{ Old = \"88\u00d717=\"; New = \"56\u00d750=\" },\n  @{ Old = \"69\u00d748=\"; New = \"36\u00d771=\" },\n  @{ Old = \"61\u00d798=\"; New = \"95\u00d771=\" },\n  @{ Old = \"90\u00d798=\"; New = \"99\u00d777=\" },\n  @{ Old = \"47\u00d792=\"; New = \"47\u00d768=\" },\n  @{ Old = \"50\u00d725=\"; New = \"28\u00d779=\" },\n  @{ Old = \"35\u00d771=\"; New = \"70\u00d741=\" },\n  @{ Old = \"73\u00d721=\"; New = \"20\u00d732=\" },\n  @{ Old = \"95\u00d738=\"; New = \"76\u00d786=\" },\n  @{ Old = \"54\u00d782=\"; New = \"56\u00d720=\" },\n  @{ Old = \"25\u00d788=\"; New = \"45\u00d740=\" },\n  @{ Old = \"52\u00d788=\"; New = \"16\u00d715=\" },\n  @{ Old = \"93\u00d732=\"; New = \"46\u00d761=\" }\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair.Old\n  $find.Replacement.Text = $pair.New\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
